$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.133.15"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 4).Value = "1.823.02"
$ws.Cells.Item(3, 5).Value = "  -0.64%  "
$ws.Cells.Item(4, 4).Value = "'1.010"
$ws.Cells.Item(4, 5).Value = "  -0.27%  "
$ws.Cells.Item(5, 4).Value = "'312.19"
$ws.Cells.Item(5, 5).Value = "  -0.57%  "
$ws.Cells.Item(6, 5).Value = "  -0.29%  "
$ws.Cells.Item(7, 4).Value = "'0.4624"
$ws.Cells.Item(7, 5).Value = "  -1.94%  "
$ws.Cells.Item(8, 4).Value = "'0.3631"
$ws.Cells.Item(8, 5).Value = "  -1.45%  "
$ws.Cells.Item(9, 4).Value = "'0.07307"
$ws.Cells.Item(9, 5).Value = "  -1.60%  "
$ws.Cells.Item(10, 4).Value = "'0.8709"
$ws.Cells.Item(10, 5).Value = "  -1.29%  "
$ws.Cells.Item(11, 4).Value = "'20.14"
$ws.Cells.Item(11, 5).Value = "  -1.52%  "
$ws.Cells.Item(12, 4).Value = "1.876.96"
$ws.Cells.Item(12, 5).Value = "  +1.61%  "
$ws.Cells.Item(13, 4).Value = "'0.07623"
$ws.Cells.Item(13, 5).Value = "  +3.89%  "
$ws.Cells.Item(14, 4).Value = "'5.341"
$ws.Cells.Item(14, 5).Value = "  -2.47%  "
$ws.Cells.Item(15, 4).Value = "'92.45"
$ws.Cells.Item(15, 5).Value = "  -0.69%  "
$ws.Cells.Item(16, 4).Value = "'6.474"
$ws.Cells.Item(16, 5).Value = "  -1.54%  "
$ws.Cells.Item(17, 5).Value = "  -0.42%  "
$ws.Cells.Item(18, 4).Value = "'0.000008653"
$ws.Cells.Item(18, 5).Value = "  -1.73%  "
$ws.Cells.Item(19, 4).Value = "'1.009"
$ws.Cells.Item(19, 5).Value = "  -0.18%  "
$ws.Cells.Item(20, 4).Value = "27.479.06"
$ws.Cells.Item(20, 5).Value = "  +0.87%  "
$ws.Cells.Item(21, 4).Value = "'14.49"
$ws.Cells.Item(21, 5).Value = "  -2.07%  "
$ws.Cells.Item(22, 4).Value = "'5.218"
$ws.Cells.Item(22, 5).Value = "  -1.63%  "
$ws.Cells.Item(23, 4).Value = "'10.56"
$ws.Cells.Item(23, 5).Value = "  -1.31%  "
$ws.Cells.Item(24, 4).Value = "2.101.13"
$ws.Cells.Item(24, 5).Value = "  +1.55%  "
$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(25, 4).Value = "'151.41"
$ws.Cells.Item(25, 5).Value = "  -1.05%  "
$ws.Cells.Item(26, 2).Value = "Toncoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).Value = "'1.860"
$ws.Cells.Item(26, 5).Value = "  -2.37%  "
$ws.Cells.Item(27, 4).Value = "'18.17"
$ws.Cells.Item(27, 5).Value = "  -2.52%  "
$ws.Cells.Item(28, 4).Value = "'2.080"
$ws.Cells.Item(28, 5).Value = "  -4.26%  "
$ws.Cells.Item(29, 4).Value = "'5.100"
$ws.Cells.Item(29, 5).Value = "  -3.42%  "
$ws.Cells.Item(30, 4).Value = "'116.07"
$ws.Cells.Item(30, 5).Value = "  -1.46%  "
$ws.Cells.Item(31, 4).Value = "'0.08900"
$ws.Cells.Item(31, 5).Value = "  -0.31%  "
$ws.Cells.Item(32, 4).Value = "'2.958"
$ws.Cells.Item(32, 5).Value = "  +0.46%  "
$ws.Cells.Item(33, 4).Value = "'0.7386"
$ws.Cells.Item(33, 5).Value = "  -2.80%  "
$ws.Cells.Item(34, 4).Value = "'4.456"
$ws.Cells.Item(34, 5).Value = "  -1.99%  "
$ws.Cells.Item(35, 4).Value = "'1.140"
$ws.Cells.Item(35, 5).Value = "  -2.82%  "
$ws.Cells.Item(36, 4).Value = "'1.010"
$ws.Cells.Item(36, 5).Value = "  -0.18%  "
$ws.Cells.Item(37, 4).Value = "'2.482"
$ws.Cells.Item(37, 5).Value = "  +3.21%  "
$ws.Cells.Item(38, 4).Value = "'1.070"
$ws.Cells.Item(38, 5).Value = "  -3.07%  "
$ws.Cells.Item(39, 4).Value = "'0.05247"
$ws.Cells.Item(39, 5).Value = "  -1.75%  "
$ws.Cells.Item(40, 4).Value = "'0.01916"
$ws.Cells.Item(40, 5).Value = "  -2.24%  "
$ws.Cells.Item(41, 5).Value = "  -2.44%  "
$ws.Cells.Item(42, 4).Value = "'7.179"
$ws.Cells.Item(42, 5).Value = "  -2.12%  "
$ws.Cells.Item(43, 4).Value = "'0.5207"
$ws.Cells.Item(43, 5).Value = "  -2.67%  "
$ws.Cells.Item(44, 4).Value = "'0.1629"
$ws.Cells.Item(44, 5).Value = "  -2.05%  "
$ws.Cells.Item(45, 4).Value = "'8.291"
$ws.Cells.Item(45, 5).Value = "  -2.94%  "
$ws.Cells.Item(46, 4).Value = "'0.4841"
$ws.Cells.Item(46, 5).Value = "  -2.40%  "
$ws.Cells.Item(47, 4).Value = "'1.009"
$ws.Cells.Item(47, 5).Value = "  -0.28%  "
$ws.Cells.Item(48, 4).Value = "'10.14"
$ws.Cells.Item(48, 5).Value = "  -3.88%  "
$ws.Cells.Item(49, 4).Value = "'103.36"
$ws.Cells.Item(49, 5).Value = "  -0.62%  "
$ws.Cells.Item(50, 4).Value = "'1.635"
$ws.Cells.Item(50, 5).Value = "  -2.34%  "
$ws.Cells.Item(51, 4).Value = "'0.06265"
$ws.Cells.Item(51, 5).Value = "  -0.83%  "
